$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
# Row 76
$ws.Range("H76").Value = 3520.35
$ws.Range("I76").Value = 3300.1538
$ws.Range("J76").Value = 3929.2856
$ws.Range("K76").Value = 3300.1538
$ws.Range("L76").Value = 3929.2856
$ws.Range("M76").Value = -2985.1538
$ws.Range("N76").Value = -4559.2856
# Row 79
$ws.Range("H79").Value = 3520.35
$ws.Range("I79").Value = 3300.1538
$ws.Range("J79").Value = 3929.2856
$ws.Range("K79").Value = 3300.1538
$ws.Range("L79").Value = 3929.2856
$ws.Range("M79").Value = -2208.1538
$ws.Range("N79").Value = -6113.2856
# Row 129
$ws.Range("H129").Value = 3165937.5
$ws.Range("J129").Value = 1449.4783
$ws.Range("L129").Value = 4348.4349
$ws.Range("N129").Value = -14348.4349
# Row 138
$ws.Range("H138").Value = 3148.9822
$ws.Range("I138").Value = 1532.16
$ws.Range("J138").Value = 4452.871
$ws.Range("K138").Value = 4596.48
$ws.Range("L138").Value = 13358.613
$ws.Range("M138").Value = 543.5199999999995
$ws.Range("N138").Value = -23638.613

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5291.315
$ws.Range("I32").Value = 3758.4575
$ws.Range("J32").Value = 11751.214
$ws.Range("K32").Value = 3758.4575
$ws.Range("L32").Value = 11751.214
$ws.Range("M32").Value = -3471.4575
$ws.Range("N32").Value = -12325.214
# Row 61
$ws.Range("H61").Value = 4099.5
$ws.Range("I61").Value = 2132.6667
$ws.Range("K61").Value = 2132.6667
$ws.Range("M61").Value = -1920.6667
# Row 74
$ws.Range("H74").Value = 1020.03845
$ws.Range("I74").Value = 909.8570999999999
$ws.Range("J74").Value = 1482.8
$ws.Range("K74").Value = 909.8570999999999
$ws.Range("L74").Value = 1482.8
$ws.Range("M74").Value = -35.85709999999995
$ws.Range("N74").Value = -3230.8
# Row 77
$ws.Range("H77").Value = 1020.03845
$ws.Range("I77").Value = 909.8570999999999
$ws.Range("J77").Value = 1482.8
$ws.Range("K77").Value = 4549.2855
$ws.Range("L77").Value = 7414
$ws.Range("M77").Value = -181.2855
$ws.Range("N77").Value = -16150
# Row 132
$ws.Range("H132").Value = 25644526
$ws.Range("I132").Value = 29414986
$ws.Range("J132").Value = 5399.8
$ws.Range("K132").Value = 88244958
$ws.Range("L132").Value = 16199.4
$ws.Range("M132").Value = -88242428
$ws.Range("N132").Value = -21259.4
# Row 136
$ws.Range("H136").Value = 4099.5
$ws.Range("I136").Value = 2132.6667
$ws.Range("K136").Value = 6398.000100000001
$ws.Range("M136").Value = -3848.000100000001
# Row 138
$ws.Range("H138").Value = 33030.645
$ws.Range("J138").Value = 33030.645
$ws.Range("L138").Value = 33030.645
$ws.Range("N138").Value = -43310.645

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 394.3158
$ws.Range("I94").Value = 371.77777
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 371.77777
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 79.22223000000002
$ws.Range("N94").Value = -1702
# Row 134
$ws.Range("H134").Value = 3946.4666
$ws.Range("I134").Value = 2433.0833
$ws.Range("K134").Value = 7299.249899999999
$ws.Range("M134").Value = -4764.249899999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 19233778
$ws.Range("I58").Value = 1694.8125
$ws.Range("J58").Value = 50005110
$ws.Range("K58").Value = 1694.8125
$ws.Range("L58").Value = 50005110
$ws.Range("M58").Value = -1491.8125
$ws.Range("N58").Value = -50005516
# Row 122
$ws.Range("H122").Value = 3206.2144
$ws.Range("J122").Value = 3677.8333
$ws.Range("L122").Value = 11033.4999
$ws.Range("N122").Value = -15933.4999
# Row 132
$ws.Range("H132").Value = 3555.8918
$ws.Range("I132").Value = 2532.6155
$ws.Range("J132").Value = 4110.1665
$ws.Range("K132").Value = 7597.8465
$ws.Range("L132").Value = 12330.4995
$ws.Range("M132").Value = -5067.8465
$ws.Range("N132").Value = -17390.4995
# Row 134
$ws.Range("H134").Value = 3710.0908
$ws.Range("I134").Value = 2601.375
$ws.Range("J134").Value = 6666.6665
$ws.Range("K134").Value = 7804.125
$ws.Range("L134").Value = 19999.9995
$ws.Range("M134").Value = -5269.125
$ws.Range("N134").Value = -25069.9995
# Row 136
$ws.Range("H136").Value = 19233778
$ws.Range("I136").Value = 1694.8125
$ws.Range("J136").Value = 50005110
$ws.Range("K136").Value = 5084.4375
$ws.Range("L136").Value = 150015330
$ws.Range("M136").Value = -2534.4375
$ws.Range("N136").Value = -150020430

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("H69").Value = 159400
$ws.Range("I69").Value = 1900
$ws.Range("J69").Value = 185650
$ws.Range("K69").Value = 5700
$ws.Range("L69").Value = 556950
$ws.Range("M69").Value = -4889
$ws.Range("N69").Value = -558572
# Row 72
$ws.Range("H72").Value = 159400
$ws.Range("I72").Value = 1900
$ws.Range("J72").Value = 185650
$ws.Range("K72").Value = 17100
$ws.Range("L72").Value = 1670850
$ws.Range("M72").Value = -13044
$ws.Range("N72").Value = -1678962
# Row 131
$ws.Range("H131").Value = 959.5806
$ws.Range("I131").Value = 816.5625
$ws.Range("J131").Value = 1112.1333
$ws.Range("K131").Value = 2449.6875
$ws.Range("L131").Value = 3336.3999
$ws.Range("M131").Value = 2590.3125
$ws.Range("N131").Value = -13416.3999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
# Row 70
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -4540
# Row 73
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -5872
# Row 97
$ws.Range("H97").Value = 1439.5
$ws.Range("I97").Value = 981.4286
$ws.Range("J97").Value = 3042.75
$ws.Range("K97").Value = 981.4286
$ws.Range("L97").Value = 3042.75
$ws.Range("M97").Value = -485.4286
$ws.Range("N97").Value = -4034.75
# Row 132
$ws.Range("H132").Value = 2523.1143
$ws.Range("I132").Value = 1773.5
$ws.Range("J132").Value = 4688.6665
$ws.Range("K132").Value = 5320.5
$ws.Range("L132").Value = 14065.9995
$ws.Range("M132").Value = -2790.5
$ws.Range("N132").Value = -19125.9995
# Row 141
$ws.Range("H141").Value = 37542.855
$ws.Range("J141").Value = 37542.855
$ws.Range("L141").Value = 37542.855
$ws.Range("N141").Value = -47902.855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 87
$ws.Range("H87").Value = 29333.334
$ws.Range("J87").Value = 29333.334
$ws.Range("L87").Value = 29333.334
$ws.Range("N87").Value = -31579.334
# Row 90
$ws.Range("H90").Value = 29333.334
$ws.Range("J90").Value = 29333.334
$ws.Range("L90").Value = 88000.00199999999
$ws.Range("N90").Value = -99232.00199999999
# Row 132
$ws.Range("H132").Value = 2457.2449
$ws.Range("I132").Value = 1665.5172
$ws.Range("J132").Value = 3605.25
$ws.Range("K132").Value = 4996.5516
$ws.Range("L132").Value = 10815.75
$ws.Range("M132").Value = -2466.5516
$ws.Range("N132").Value = -15875.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16040
# Row 132
$ws.Range("H132").Value = 179446.61
$ws.Range("I132").Value = 234476.25
$ws.Range("J132").Value = 10427
$ws.Range("K132").Value = 703428.75
$ws.Range("L132").Value = 31281
$ws.Range("M132").Value = -700898.75
$ws.Range("N132").Value = -36341

